$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster data (Player, Position, Team) for rows 2-19
$data = @(
    @("Immanuel Quickley", "PG,SG", "Toronto Raptors"),
    @("Devin Booker", "PG,SG", "Phoenix Suns"),
    @("Trae Young", "PG", "Atlanta Hawks"),
    @("LeBron James", "SF,PF", "Los Angeles Lakers"),
    @("P.J. Washington", "PF", "Dallas Mavericks"),
    @("Gradey Dick", "SG,SF", "Toronto Raptors"),
    @("Desmond Bane", "SG,SF", "Memphis Grizzlies"),
    @("Walker Kessler", "C", "Utah Jazz"),
    @("Norman Powell", "SG,SF", "LA Clippers"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Coby White", "PG,SG", "Chicago Bulls"),
    @("Jalen Brunson", "PG", "New York Knicks"),
    @("Devin Vassell", "SG,SF", "San Antonio Spurs"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("D'Angelo Russell", "PG", "Brooklyn Nets"),
    @("Brandon Ingram", "SG,SF,PF", "New Orleans Pelicans"),
    @("Jabari Smith Jr.", "PF,C", "Houston Rockets"),
    @("Kawhi Leonard", "SG,SF,PF", "LA Clippers")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
